# Insert a new data row at row 555 (pushing the existing rows 555..622 down
# to 556..623) and populate it with the new record's values. All the other
# "static" columns in this table follow a fixed pattern shared by every row
# (Mercado ID, Mercado, Region, Codreg, Categoria ID, Categoria, Variedad,
# Calidad, Unidad de comercializacion, Origen, Kg o Unidades, Clasificacion),
# so we copy that pattern from the row immediately below (the row that used
# to be 555, now pushed to 556) to keep everything consistent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 555:622 down to 556:623, creating a blank row 555.
$ws.Rows.Item(555).Insert()

# New row's own data.
$newDate = 45212
$newVolumen = 220
$newPrecioMin = 3800
$newPrecioMax = 4000
$newPrecioProm = 3900
$newPrecioKg = 1300

# Static / repeated columns, copied from the pattern used by every row in
# this table.
$ws.Cells.Item(555, 1).Value = 3
$ws.Cells.Item(555, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(555, 3).Value = "Coquimbo"
$ws.Cells.Item(555, 4).Value = $newDate
$ws.Cells.Item(555, 4).NumberFormat = $ws.Cells.Item(556, 4).NumberFormat
$ws.Cells.Item(555, 5).Value = 5
$ws.Cells.Item(555, 6).Value = 100112012
$ws.Cells.Item(555, 7).Value = "Espinaca"
$ws.Cells.Item(555, 8).Value = "Sin especificar"
$ws.Cells.Item(555, 9).Value = "Primera"
$ws.Cells.Item(555, 10).Value = $newVolumen
$ws.Cells.Item(555, 11).Value = $newPrecioMin
$ws.Cells.Item(555, 12).Value = $newPrecioMax
$ws.Cells.Item(555, 13).Value = $newPrecioProm
$ws.Cells.Item(555, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(555, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(555, 16).Value = $newPrecioKg
$ws.Cells.Item(555, 17).Value = 3
$ws.Cells.Item(555, 18).Value = "Hortaliza"
